$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fill in the new risk row (row 9) with the new data.
$ws.Range("A9").Value = "La no aseptacion de la propuesta de aceptacion de elaboracion del producto por parte  de cliente."
$ws.Range("B9").Value = "Replanteamiento de la propuesta según criticas del usuario."
$ws.Range("C9").Value = "MEDIA"
$ws.Range("D9").Value = "Observacion de cada uno de los puntos dados por el usuario y generalizacion entre el equipo de los requerimientos del cliente."

# 2. Apply a full grid (inside borders) across the whole table range, including
#    the two new blank trailing rows (10-11) that close the table visually.
$rng = $ws.Range("A1:E11")
$rng.Borders.Item(11).LineStyle = 1
$rng.Borders.Item(12).LineStyle = 1

# 3. E9 was filled by copying B9 (same "PLAN MITIGACION" text ended up in the
#    "PLAN CONTINGENCIA" column) which also carries over B9's border style.
$ws.Range("B9").Copy()
$ws.Range("E9").PasteSpecial()

# 4. Row 10 gets a distinguishing underline font on A10 (cursor left there).
$ws.Range("A10").Font.Underline = 2

# 5. Row heights.
$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 23.25

# 6. Column widths (manual resize, no longer auto "best fit").
$ws.Columns.Item(1).ColumnWidth = 50.83
$ws.Columns.Item(2).ColumnWidth = 47.0
$ws.Columns.Item(3).ColumnWidth = 7.67
$ws.Columns.Item(4).ColumnWidth = 58.0

# 7. Move the active selection to A10, matching where the user left off.
$ws.Range("A10").Select()
